$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 73.53729552766163
$ws.Range("B3").Value = 0.9167350811896852
$ws.Range("B4").Value = 0.05918277294611446
$ws.Range("B5").Value = 0.373696506396882
